$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 2327.3845
$ws.Cells.Item(58, 10).Value = 3241.5
$ws.Cells.Item(58, 12).Value = 9724.5
$ws.Cells.Item(58, 14).Value = -10024.5

$ws.Cells.Item(68, 8).Value = 50295
$ws.Cells.Item(68, 10).Value = 50295
$ws.Cells.Item(68, 12).Value = 50295
$ws.Cells.Item(68, 14).Value = -51793

$ws.Cells.Item(70, 8).Value = 4502.08
$ws.Cells.Item(70, 10).Value = 5661.2144
$ws.Cells.Item(70, 12).Value = 16983.6432
$ws.Cells.Item(70, 14).Value = -17523.6432

$ws.Cells.Item(71, 8).Value = 50295
$ws.Cells.Item(71, 10).Value = 50295
$ws.Cells.Item(71, 12).Value = 150885
$ws.Cells.Item(71, 14).Value = -158373

$ws.Cells.Item(73, 8).Value = 4502.08
$ws.Cells.Item(73, 10).Value = 5661.2144
$ws.Cells.Item(73, 12).Value = 16983.6432
$ws.Cells.Item(73, 14).Value = -18855.6432

$ws.Cells.Item(86, 8).Value = 2875.625
$ws.Cells.Item(86, 10).Value = 3000.6667
$ws.Cells.Item(86, 12).Value = 3000.6667
$ws.Cells.Item(86, 14).Value = -5246.6667

$ws.Cells.Item(89, 8).Value = 2875.625
$ws.Cells.Item(89, 10).Value = 3000.6667
$ws.Cells.Item(89, 12).Value = 15003.3335
$ws.Cells.Item(89, 14).Value = -26235.3335

$ws.Cells.Item(92, 8).Value = 271
$ws.Cells.Item(92, 9).Value = 233.33333
$ws.Cells.Item(92, 11).Value = 233.33333
$ws.Cells.Item(92, 13).Value = 1014.66667

$ws.Cells.Item(116, 8).Value = 3992.5
$ws.Cells.Item(116, 10).Value = 3992.5
$ws.Cells.Item(116, 12).Value = 3992.5
$ws.Cells.Item(116, 14).Value = -10876.5

$ws.Cells.Item(129, 8).Value = 1827
$ws.Cells.Item(129, 9).Value = 659
$ws.Cells.Item(129, 11).Value = 1977
$ws.Cells.Item(129, 13).Value = 3023

$ws.Cells.Item(135, 8).Value = 761.9286
$ws.Cells.Item(135, 9).Value = 615.1667
$ws.Cells.Item(135, 11).Value = 5536.5003
$ws.Cells.Item(135, 13).Value = -3001.5003

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2056.5789
$ws.Cells.Item(61, 9).Value = 1370.7333
$ws.Cells.Item(61, 10).Value = 4628.5
$ws.Cells.Item(61, 11).Value = 1370.7333
$ws.Cells.Item(61, 12).Value = 4628.5
$ws.Cells.Item(61, 13).Value = -1158.7333
$ws.Cells.Item(61, 14).Value = -5052.5

$ws.Cells.Item(74, 8).Value = 6062
$ws.Cells.Item(74, 9).Value = 5545.4287
$ws.Cells.Item(74, 11).Value = 5545.4287
$ws.Cells.Item(74, 13).Value = -4671.4287

$ws.Cells.Item(77, 8).Value = 6062
$ws.Cells.Item(77, 9).Value = 5545.4287
$ws.Cells.Item(77, 11).Value = 27727.1435
$ws.Cells.Item(77, 13).Value = -23359.1435

$ws.Cells.Item(97, 8).Value = 1346.2727
$ws.Cells.Item(97, 9).Value = 981
$ws.Cells.Item(97, 11).Value = 981
$ws.Cells.Item(97, 13).Value = -485

$ws.Cells.Item(110, 8).Value = 100002110
$ws.Cells.Item(110, 9).Value = 125001830
$ws.Cells.Item(110, 11).Value = 125001830
$ws.Cells.Item(110, 13).Value = -124999785

$ws.Cells.Item(136, 8).Value = 2056.5789
$ws.Cells.Item(136, 9).Value = 1370.7333
$ws.Cells.Item(136, 10).Value = 4628.5
$ws.Cells.Item(136, 11).Value = 4112.199900000001
$ws.Cells.Item(136, 12).Value = 13885.5
$ws.Cells.Item(136, 13).Value = -1562.199900000001
$ws.Cells.Item(136, 14).Value = -18985.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3582.3333
$ws.Cells.Item(20, 9).Value = 3123.75
$ws.Cells.Item(20, 10).Value = 4499.5
$ws.Cells.Item(20, 11).Value = 3123.75
$ws.Cells.Item(20, 12).Value = 4499.5
$ws.Cells.Item(20, 13).Value = -2876.75
$ws.Cells.Item(20, 14).Value = -4993.5

$ws.Cells.Item(36, 8).Value = 1481.75
$ws.Cells.Item(36, 9).Value = 1481.75
$ws.Cells.Item(36, 11).Value = 1481.75
$ws.Cells.Item(36, 13).Value = -947.75

$ws.Cells.Item(86, 8).Value = 4930.0713
$ws.Cells.Item(86, 9).Value = 1655.6666
$ws.Cells.Item(86, 11).Value = 1655.6666
$ws.Cells.Item(86, 13).Value = -532.6666

$ws.Cells.Item(89, 8).Value = 4930.0713
$ws.Cells.Item(89, 9).Value = 1655.6666
$ws.Cells.Item(89, 11).Value = 8278.333000000001
$ws.Cells.Item(89, 13).Value = -2662.333000000001

$ws.Cells.Item(99, 8).Value = 250000700
$ws.Cells.Item(99, 9).Value = 333333920
$ws.Cells.Item(99, 11).Value = 333333920
$ws.Cells.Item(99, 13).Value = -333332422

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(68, 8).Value = 88382.5
$ws.Cells.Item(68, 10).Value = 88382.5
$ws.Cells.Item(68, 12).Value = 88382.5
$ws.Cells.Item(68, 14).Value = -89880.5

$ws.Cells.Item(71, 8).Value = 88382.5
$ws.Cells.Item(71, 10).Value = 88382.5
$ws.Cells.Item(71, 12).Value = 265147.5
$ws.Cells.Item(71, 14).Value = -272635.5

$ws.Cells.Item(107, 8).Value = 1522.9286
$ws.Cells.Item(107, 10).Value = 2775.5
$ws.Cells.Item(107, 12).Value = 2775.5
$ws.Cells.Item(107, 14).Value = -6615.5

$ws.Cells.Item(138, 8).Value = 130000
$ws.Cells.Item(138, 10).Value = 130000
$ws.Cells.Item(138, 12).Value = 130000
$ws.Cells.Item(138, 14).Value = -140280

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 1000
$ws.Cells.Item(51, 10).Value = 1500
$ws.Cells.Item(51, 12).Value = 4500
$ws.Cells.Item(51, 14).Value = -5420

$ws.Cells.Item(59, 8).Value = 900
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).ClearContents()

$ws.Cells.Item(80, 8).Value = 4291.7393
$ws.Cells.Item(80, 9).Value = 4021.611
$ws.Cells.Item(80, 11).Value = 12064.833
$ws.Cells.Item(80, 13).Value = -11128.833

$ws.Cells.Item(83, 8).Value = 4291.7393
$ws.Cells.Item(83, 9).Value = 4021.611
$ws.Cells.Item(83, 11).Value = 36194.499
$ws.Cells.Item(83, 13).Value = -31514.499

$ws.Cells.Item(140, 8).Value = 2383.5625
$ws.Cells.Item(140, 9).Value = 1654
$ws.Cells.Item(140, 11).Value = 4962
$ws.Cells.Item(140, 13).Value = 218

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4490
$ws.Cells.Item(80, 9).Value = 4490
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 4490
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = -3492
$ws.Cells.Item(80, 14).ClearContents()

$ws.Cells.Item(83, 8).Value = 4490
$ws.Cells.Item(83, 9).Value = 4490
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 22450
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = -17458
$ws.Cells.Item(83, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 2720.8572
$ws.Cells.Item(93, 9).Value = 2832.6667
$ws.Cells.Item(93, 11).Value = 2832.6667
$ws.Cells.Item(93, 13).Value = -1584.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value = 3000
$ws.Cells.Item(12, 9).Value = 3000
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 3000
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = -2858
$ws.Cells.Item(12, 14).ClearContents()

$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).ClearContents()
